# HeatMap para análise inicial de dados
# -------------------------------------------------------------
# Summary of changes applied:
#  1. Insert a new worksheet "Plan1" (a small binary-to-decimal
#     "heat map" helper table) right before "Correlação".
#  2. Update the saved cursor/selection on a few existing sheets.
#  3. Tweak a handful of row heights on "Correlação" (some grow to
#     fit wrapped text, others collapse back to the sheet default).
# -------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. "Colunas e Relação com Vazio" : just move the saved selection ---
$wsColunas = $wb.Worksheets.Item("Colunas e Relação com Vazio")
$wsColunas.Activate()
$wsColunas.Range("B4").Select()

# --- 2. "Correlação" : move selection + adjust a batch of row heights ---
$wsCorr = $wb.Worksheets.Item("Correlação")
$wsCorr.Activate()
$wsCorr.Range("B4").Select()

# Rows whose text now wraps into two (45pt) or three (60pt) lines.
$rowsTo60 = @(28, 39)
$rowsTo45 = @(29, 30, 40, 44, 53, 55, 56, 57, 58)
# Rows that collapse back down to the sheet's default height.
$rowsAutoFit = @(31, 32, 34, 36, 37, 38, 41, 42, 45, 46, 47, 48, 49, 52, 54)

foreach ($r in $rowsTo60) {
    $wsCorr.Rows.Item($r).RowHeight = 60
}
foreach ($r in $rowsTo45) {
    $wsCorr.Rows.Item($r).RowHeight = 45
}
foreach ($r in $rowsAutoFit) {
    $wsCorr.Rows.Item($r).AutoFit()
}

# --- 3. "Plan2" : move the saved selection ---
$wsPlan2 = $wb.Worksheets.Item("Plan2")
$wsPlan2.Activate()
$wsPlan2.Range("A26:B29").Select()

# --- 4. Insert the new "Plan1" sheet right before "Correlação" ---
$beforeSheet = $wb.Worksheets.Item("Correlação")
$plan1 = $wb.Worksheets.Add($beforeSheet)
$plan1.Name = "Plan1"

# Row 1: bit weights (3,2,1,0)
$plan1.Range("B1").Value = 3
$plan1.Range("C1").Value = 2
$plan1.Range("D1").Value = 1
$plan1.Range("E1").Value = 0

# Row 2: powers of two for each weight (B2 alone, C2:E2 filled as one shared formula)
$plan1.Range("B2").Formula = "=2^B1"
$plan1.Range("C2:E2").Formula = "=2^C1"

# Row 3: column headers (existing shared strings A1_4 / A1_3 / A1_2 / A1_1)
$plan1.Range("B3").Value = "A1_4"
$plan1.Range("C3").Value = "A1_3"
$plan1.Range("D3").Value = "A1_2"
$plan1.Range("E3").Value = "A1_1"

# Row 4: all-zero combination, with its own (non-shared) total formula
$plan1.Range("B4").Value = 0
$plan1.Range("C4").Value = 0
$plan1.Range("D4").Value = 0
$plan1.Range("E4").Value = 0
$plan1.Range("F4").Formula = "=B4*`$B`$2+C4*`$C`$2+D4*`$D`$2+E4*`$E`$2"

# Rows 5-11: the rest of the 3-bit binary combinations (B column left blank)
$plan1.Range("C5").Value = 0
$plan1.Range("D5").Value = 0
$plan1.Range("E5").Value = 1

$plan1.Range("C6").Value = 0
$plan1.Range("D6").Value = 1
$plan1.Range("E6").Value = 0

$plan1.Range("C7").Value = 1
$plan1.Range("D7").Value = 0
$plan1.Range("E7").Value = 0

$plan1.Range("C8").Value = 0
$plan1.Range("D8").Value = 1
$plan1.Range("E8").Value = 1

$plan1.Range("C9").Value = 1
$plan1.Range("D9").Value = 0
$plan1.Range("E9").Value = 1

$plan1.Range("C10").Value = 1
$plan1.Range("D10").Value = 1
$plan1.Range("E10").Value = 0

$plan1.Range("C11").Value = 1
$plan1.Range("D11").Value = 1
$plan1.Range("E11").Value = 1

# Column F (rows 5-17): shared "weighted total" formula, filled down.
$plan1.Range("F5:F17").Formula = "=B5*`$B`$2+C5*`$C`$2+D5*`$D`$2+E5*`$E`$2"

# Saved selection / active cell on the new sheet.
$plan1.Range("C3:F11").Select()

Write-Host "Plan1 sheet created and existing sheets updated."
